# "New crime data collected" -- weekly CompStat refresh for the 111th Precinct.
# Updates the report header (volume/number + week-covering dates) and
# refreshes every statistic in the Crime Complaints table (rows 15-27)
# to reflect the newly collected week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (rich-text shared strings) -- edit the runs in place so that
# formatting/other runs are preserved.
# ---------------------------------------------------------------------------

# A8: "Volume 30   Number  43" -> "...  44"
$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$numStart = $volText.Length - 1
$volCell.Characters($numStart, 2).Text = "44"

# C9: "Report Covering the Week  10/23/2023  Through  10/29/2023"
#     -> "...10/30/2023  Through  11/5/2023"
$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 10).Text = "10/30/2023"
$weekCell.Characters(48, 10).Text = "11/5/2023"

# ---------------------------------------------------------------------------
# Crime Complaints table -- numeric refresh.
# Row map: C/D = Week to Date 2023/2022, E = %Chg, F/G = 28 Day 2023/2022,
#          H = %Chg, I/J = YTD 2023/2022, K = %Chg, L = 2Yr %Chg,
#          M = 13Yr %Chg, N = 30Yr %Chg.
# ---------------------------------------------------------------------------

# Row 15 (Rape)
$ws.Range("N15").Value2 = -33.333333333333

# Row 16 (Robbery)
$ws.Range("D16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 8
$ws.Range("H16").Value2 = -87.5
$ws.Range("J16").Value2 = 73
$ws.Range("K16").Value2 = -30.136986301369
$ws.Range("L16").Value2 = 131.818181818182
$ws.Range("M16").Value2 = -25
$ws.Range("N16").Value2 = -81.386861313868

# Row 17 (Fel. Assault)
$ws.Range("C17").Value2 = 3
$ws.Range("D17").Value2 = 7
$ws.Range("E17").Value2 = -57.142857142857
$ws.Range("F17").Value2 = 11
$ws.Range("G17").Value2 = 16
$ws.Range("H17").Value2 = -31.25
$ws.Range("I17").Value2 = 83
$ws.Range("J17").Value2 = 86
$ws.Range("K17").Value2 = -3.488372093023
$ws.Range("L17").Value2 = 31.746031746031
$ws.Range("M17").Value2 = 80.434782608695
$ws.Range("N17").Value2 = -22.429906542056

# Row 18 (Burglary)
$ws.Range("D18").Value2 = 7
$ws.Range("E18").Value2 = 28.571428571428
$ws.Range("F18").Value2 = 29
$ws.Range("G18").Value2 = 23
$ws.Range("H18").Value2 = 26.086956521739
$ws.Range("I18").Value2 = 264
$ws.Range("J18").Value2 = 229
$ws.Range("K18").Value2 = 15.283842794759
$ws.Range("L18").Value2 = 42.702702702702
$ws.Range("M18").Value2 = 26.315789473684
$ws.Range("N18").Value2 = -69.515011547344

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value2 = 16
$ws.Range("D19").Value2 = 21
$ws.Range("E19").Value2 = -23.809523809523
$ws.Range("F19").Value2 = 52
$ws.Range("G19").Value2 = 69
$ws.Range("H19").Value2 = -24.637681159420
$ws.Range("I19").Value2 = 540
$ws.Range("J19").Value2 = 563
$ws.Range("K19").Value2 = -4.085257548845
$ws.Range("L19").Value2 = 70.347003154574
$ws.Range("M19").Value2 = 74.193548387096
$ws.Range("N19").Value2 = 10.429447852760

# Row 20 (G.L.A.) -- C20/D20/E20 go from text placeholders ("0" / "***.*")
# to real numbers now that last year's week-to-date count is non-zero, so
# the cell styles switch from the text style (14) to the numeric ones.
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 5
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value2 = -80
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F20").Value2 = 13
$ws.Range("G20").Value2 = 10
$ws.Range("H20").Value2 = 30
$ws.Range("I20").Value2 = 161
$ws.Range("J20").Value2 = 98
$ws.Range("K20").Value2 = 64.285714285714
$ws.Range("L20").Value2 = 182.456140350877
$ws.Range("M20").Value2 = 34.166666666666
$ws.Range("N20").Value2 = -94.260249554367

# Row 21 (TOTAL)
$ws.Range("D21").Value2 = 43
$ws.Range("E21").Value2 = -32.558139534883
$ws.Range("F21").Value2 = 106
$ws.Range("G21").Value2 = 126
$ws.Range("H21").Value2 = -15.873015873015
$ws.Range("I21").Value2 = 1109
$ws.Range("J21").Value2 = 1051
$ws.Range("K21").Value2 = 5.518553758325
$ws.Range("L21").Value2 = 70.878274268104
$ws.Range("M21").Value2 = 46.306068601583
$ws.Range("N21").Value2 = -75.658472344161

# Row 24 (Petit Larceny)
$ws.Range("C24").Value2 = 5
$ws.Range("D24").Value2 = 10
$ws.Range("E24").Value2 = -50
$ws.Range("F24").Value2 = 42
$ws.Range("G24").Value2 = 52
$ws.Range("H24").Value2 = -19.230769230769
$ws.Range("I24").Value2 = 493
$ws.Range("J24").Value2 = 659
$ws.Range("K24").Value2 = -25.189681335356
$ws.Range("L24").Value2 = 6.021505376344
$ws.Range("M24").Value2 = 30.423280423280

# Row 25 (Misd. Assault)
$ws.Range("D25").Value2 = 5
$ws.Range("E25").Value2 = -40
$ws.Range("F25").Value2 = 15
$ws.Range("G25").Value2 = 35
$ws.Range("H25").Value2 = -57.142857142857
$ws.Range("I25").Value2 = 194
$ws.Range("J25").Value2 = 211
$ws.Range("K25").Value2 = -8.056872037914
$ws.Range("L25").Value2 = 50.387596899224
$ws.Range("M25").Value2 = 30.201342281879

# Row 27 (Other Sex Crimes)
$ws.Range("L27").Value2 = -5.882352941176
